$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "First Try"
$co = $ws1.ChartObjects(1)
$chart = $co.Chart
$sc = $chart.SeriesCollection()
Write-Host "Count before: $($sc.Count)"
$old = $sc.Item(1)
$old.Delete()
Write-Host "Count after delete: $($sc.Count)"
$ns = $sc.NewSeries()
$ns.Name = "The Cube"
$ns.XValues = $ws1.Range("B17:Y17")
$ns.Values = $ws1.Range("B18:Y18")
Write-Host "Count after add: $($sc.Count)"
